# Applies the numeric "want-to-go" count refresh (and one sold-out -> on-sale
# status change) described in the commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - column F ("想去人数") updates
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 530
$ws1.Range("F5").Value  = 2536
$ws1.Range("F7").Value  = 90
$ws1.Range("F9").Value  = 1717
$ws1.Range("F10").Value = 1717
$ws1.Range("F11").Value = 1426
$ws1.Range("F12").Value = 84
$ws1.Range("F16").Value = 1003
$ws1.Range("F17").Value = 346
$ws1.Range("F20").Value = 7624
$ws1.Range("F21").Value = 8718
$ws1.Range("F31").Value = 365
$ws1.Range("F32").Value = 1535
$ws1.Range("F35").Value = 255
$ws1.Range("F41").Value = 1382
$ws1.Range("F44").Value = 230
$ws1.Range("F46").Value = 227
$ws1.Range("F49").Value = 53

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - row 18 ticket went from sold-out to on-sale
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 0
$ws2.Range("G18").Value = 80
$ws2.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=92468"
$ws2.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202409/hLIqUa8w1726239513479.png"

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - column F updates
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 197
$ws3.Range("F3").Value = 2670
$ws3.Range("F6").Value = 27

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - column F updates
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 197
$ws4.Range("F8").Value  = 530
$ws4.Range("F9").Value  = 2536
$ws4.Range("F11").Value = 90
$ws4.Range("F13").Value = 1717
$ws4.Range("F14").Value = 1717
$ws4.Range("F15").Value = 84
$ws4.Range("F18").Value = 1003
$ws4.Range("F19").Value = 346
$ws4.Range("F24").Value = 7624
$ws4.Range("F25").Value = 7624
$ws4.Range("F26").Value = 8718
$ws4.Range("F33").Value = 365
$ws4.Range("F34").Value = 1535
$ws4.Range("F37").Value = 255
$ws4.Range("F45").Value = 230
$ws4.Range("F47").Value = 227
$ws4.Range("F51").Value = 53
